$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Adenocarcinoma']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the B2 query text (Primary diagnosis query fixed)
$ws.Range("B2").Value = $newQuery

# Preserve the original row height for row 2 (unchanged by the content edit)
$ws.Rows.Item(2).RowHeight = 157.5

# Update the active selection to C2 (as recorded in the saved view)
[void]$ws.Range("C2").Select()
